# Applies the cryptos-list data refresh described in the commit:
# "Updated cryptos list on Thu Jan 25 21:18:08 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '39.859.39'
# Row 3
$ws.Range('D3').Value = '2.226.29'
$ws.Range('E3').Value = '  +0.78%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '291.97'
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.26'
$ws.Range('E6').Value = '  +1.20%  '
# Row 7
$ws.Range('E7').Value = '  -0.14%  '
# Row 8
$ws.Range('E8').Value = '  +0.00%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.468'
$ws.Range('E9').Value = '  -0.31%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.57'
$ws.Range('E10').Value = '  +1.03%  '
# Row 11
$ws.Range('E11').Value = '  -0.34%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.87'
$ws.Range('E12').Value = '  +5.09%  '
# Row 13
$ws.Range('E13').Value = '  +2.71%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.46'
$ws.Range('E14').Value = '  +2.34%  '
# Row 15
$ws.Range('D15').Value = '2.570.37'
$ws.Range('E15').Value = '  +0.83%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.85'
$ws.Range('E16').Value = '  -1.14%  '
# Row 17
$ws.Range('D17').Value = '2.236.38'
$ws.Range('E17').Value = '  +0.87%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.733'
$ws.Range('E18').Value = '  +0.88%  '
# Row 19
$ws.Range('D19').Value = '39.820.28'
$ws.Range('E19').Value = '  +0.43%  '
# Row 20
$ws.Range('D20').Value = '0.0₃0888'
$ws.Range('E20').Value = '  +1.17%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.11'
$ws.Range('E21').Value = '  -2.09%  '
# Row 22
$ws.Range('E22').Value = '  -0.34%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.79'
$ws.Range('E23').Value = '  +0.25%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '237.35'
$ws.Range('E24').Value = '  +0.83%  '
# Row 25
$ws.Range('E25').Value = '  -0.11%  '
# Row 26
$ws.Range('E26').Value = '  -0.03%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.83'
$ws.Range('E27').Value = '  +0.26%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.12'
$ws.Range('E28').Value = '  +2.07%  '
# Row 29
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.26'
$ws.Range('E29').Value = '  +0.25%  '
# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.05'
$ws.Range('E30').Value = '  -6.69%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.03'
$ws.Range('E31').Value = '  +3.55%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.02'
$ws.Range('E32').Value = '  -1.76%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.02%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.98'
$ws.Range('E34').Value = '  +1.40%  '
# Row 35
$ws.Range('E35').Value = '  +7.61%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0714'
$ws.Range('E36').Value = '  -0.12%  '
# Row 38
$ws.Range('E38').Value = '  -0.05%  '
# Row 39
$ws.Range('E39').Value = '  +4.03%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0989'
$ws.Range('E40').Value = '  +0.53%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.35'
$ws.Range('E41').Value = '  -3.76%  '
# Row 42
$ws.Range('D42').Value = '2.109.63'
$ws.Range('E42').Value = '  +2.11%  '
# Row 43
$ws.Range('E43').Value = '  -1.14%  '
# Row 44
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0272'
$ws.Range('E44').Value = '  +1.87%  '
# Row 45
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.14'
$ws.Range('E45').Value = '  +2.29%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.00'
$ws.Range('E46').Value = '  +0.54%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.00'
$ws.Range('E47').Value = '  -8.01%  '
# Row 48
$ws.Range('E48').Value = '  +4.94%  '
# Row 49
$ws.Range('D49').Value = '2.436.38'
$ws.Range('E49').Value = '  +0.74%  '
# Row 50
$ws.Range('E50').Value = '  +2.75%  '
# Row 51
$ws.Range('E51').Value = '  +2.93%  '
